$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove all existing hyperlinks before restructuring the table
$ws.Hyperlinks.Delete() | Out-Null

# Drop the rows that no longer exist in the refreshed course listing
$ws.Rows("22:31").Delete() | Out-Null

# Write the refreshed course listing data (rows 2-21)
$ws.Range("A2").Value2 = 'AZ-900 Azure Fundamentals - Microsoft Azure Fundamentals'
$ws.Range("B2").Value2 = 'IT & Software > IT Certifications'
$ws.Range("C2").Value2 = 'https://www.real.discount/offer/az-900-azure-fundamentals-microsoft-azure-fundamentals-15711'
$ws.Range("D2").Value2 = '9 hours ago'

$ws.Range("A3").Value2 = 'Web Applications Step by Step Guide Part - 3'
$ws.Range("B3").Value2 = 'IT & Software > Other IT & Software'
$ws.Range("C3").Value2 = 'https://www.real.discount/offer/web-applications-step-by-step-guide-part-3-27522'
$ws.Range("D3").Value2 = '9 hours ago'

$ws.Range("A4").Value2 = 'Web Application: Step by Step Guide'
$ws.Range("B4").Value2 = 'IT & Software > Other IT & Software'
$ws.Range("C4").Value2 = 'https://www.real.discount/offer/web-application-step-by-step-guide-27524'
$ws.Range("D4").Value2 = '9 hours ago'

$ws.Range("A5").Value2 = 'Web Applications Step by Step Guide Part-2'
$ws.Range("B5").Value2 = 'IT & Software > Other IT & Software'
$ws.Range("C5").Value2 = 'https://www.real.discount/offer/web-applications-step-by-step-guide-part-2-27523'
$ws.Range("D5").Value2 = '9 hours ago'

$ws.Range("A6").Value2 = 'Web Applications Step by Step Guide Part 4'
$ws.Range("B6").Value2 = 'IT & Software > Other IT & Software'
$ws.Range("C6").Value2 = 'https://www.real.discount/offer/web-applications-step-by-step-guide-part-4-27521'
$ws.Range("D6").Value2 = '9 hours ago'

$ws.Range("A7").Value2 = '18 Crucial Cyber Security Tips'
$ws.Range("B7").Value2 = 'IT & Software > Network & Security'
$ws.Range("C7").Value2 = 'https://www.real.discount/offer/18-crucial-cyber-security-tips-29894'
$ws.Range("D7").Value2 = '9 hours ago'

$ws.Range("A8").Value2 = 'Linux Mastery: CLI & Kali Commands Practice Tests 2024 pro'
$ws.Range("B8").Value2 = 'IT & Software > IT Certifications'
$ws.Range("C8").Value2 = 'https://www.real.discount/offer/linux-mastery-cli-kali-commands-practice-tests-2024-pro-38400'
$ws.Range("D8").Value2 = '7 hours ago'

$ws.Range("A9").Value2 = 'Exam MS-900: Microsoft 365 Fundamentals Mock Exams'
$ws.Range("B9").Value2 = 'IT & Software > IT Certifications'
$ws.Range("C9").Value2 = 'https://www.real.discount/offer/exam-ms-900-microsoft-365-fundamentals-mock-exams-35283'
$ws.Range("D9").Value2 = '7 hours ago'

$ws.Range("A10").Value2 = 'Learn Embarcadero Borland C++ Builder in 1 hour'
$ws.Range("B10").Value2 = 'IT & Software > Other IT & Software'
$ws.Range("C10").Value2 = 'https://www.real.discount/offer/learn-embarcadero-borland-c-builder-in-1-hour-2'
$ws.Range("D10").Value2 = '7 hours ago'

$ws.Range("A11").Value2 = 'CompTIA Security+ (SY0-701) Practice Tests'
$ws.Range("B11").Value2 = 'IT & Software > IT Certifications'
$ws.Range("C11").Value2 = 'https://www.real.discount/offer/comptia-security-sy0-701-practice-tests-36461'
$ws.Range("D11").Value2 = '6 hours ago'

$ws.Range("A12").Value2 = '(ISC)2 Certified in Cybersecurity (CC) Practice Exams'
$ws.Range("B12").Value2 = 'IT & Software > IT Certifications'
$ws.Range("C12").Value2 = 'https://www.real.discount/offer/isc-2-certified-in-cybersecurity-cc-practice-exams-36420'
$ws.Range("D12").Value2 = '6 hours ago'

$ws.Range("A13").Value2 = '(ISC)2 Certified in Cybersecurity (CC) Practice Exams: Set 2'
$ws.Range("B13").Value2 = 'IT & Software > IT Certifications'
$ws.Range("C13").Value2 = 'https://www.real.discount/offer/isc-2-certified-in-cybersecurity-cc-practice-exams-set-2-36444'
$ws.Range("D13").Value2 = '6 hours ago'

$ws.Range("A14").Value2 = 'PCEP (30-02) Practice Exams'
$ws.Range("B14").Value2 = 'IT & Software > IT Certifications'
$ws.Range("C14").Value2 = 'https://www.real.discount/offer/pcep-30-02-practice-exams-36925'
$ws.Range("D14").Value2 = '6 hours ago'

$ws.Range("A15").Value2 = 'CSS And Javascript Crash Course'
$ws.Range("B15").Value2 = 'IT & Software > IT Certifications'
$ws.Range("C15").Value2 = 'https://www.real.discount/offer/css-and-javascript-crash-course-13022'
$ws.Range("D15").Value2 = '4 hours ago'

$ws.Range("A16").Value2 = 'ECCouncil: Certified Cybersecurity Technician'
$ws.Range("B16").Value2 = 'IT & Software > IT Certifications'
$ws.Range("C16").Value2 = 'https://www.real.discount/offer/eccouncil-certified-cybersecurity-technician-39209'
$ws.Range("D16").Value2 = '3 minutes ago'

$ws.Range("A17").Value2 = 'Web Automation and Scraping using Python'
$ws.Range("B17").Value2 = 'IT & Software > Other IT & Software'
$ws.Range("C17").Value2 = 'https://www.real.discount/offer/web-automation-and-scraping-using-python-36771'
$ws.Range("D17").Value2 = '12 hours ago'

$ws.Range("A18").Value2 = 'The Best ChatGPT & AI Course: Make Money With AI'
$ws.Range("B18").Value2 = 'IT & Software > Other IT & Software'
$ws.Range("C18").Value2 = 'https://www.real.discount/offer/the-best-chatgpt-ai-course-make-money-with-ai-35563'
$ws.Range("D18").Value2 = '12 hours ago'

$ws.Range("A19").Value2 = 'Learn Azure Bicep'
$ws.Range("B19").Value2 = 'IT & Software > Other IT & Software'
$ws.Range("C19").Value2 = 'https://www.real.discount/offer/learn-azure-bicep-20512'
$ws.Range("D19").Value2 = '11 hours ago'

$ws.Range("A20").Value2 = 'Midjourney for Beginners: Embark on Your Artistic Journey'
$ws.Range("B20").Value2 = 'IT & Software > Other IT & Software'
$ws.Range("C20").Value2 = 'https://www.real.discount/offer/midjourney-for-beginners-embark-on-your-artistic-journey-35457'
$ws.Range("D20").Value2 = '11 hours ago'

$ws.Range("A21").Value2 = 'Google Forms o Formularios de Cero a Avanzado'
$ws.Range("B21").Value2 = 'IT & Software > Other IT & Software'
$ws.Range("C21").Value2 = 'https://www.real.discount/offer/google-forms-o-formularios-de-cero-a-avanzado-32061'
$ws.Range("D21").Value2 = '11 hours ago'

# Re-create hyperlinks on the Link column for each data row
$ws.Hyperlinks.Add($ws.Range("C2"), 'https://www.real.discount/offer/az-900-azure-fundamentals-microsoft-azure-fundamentals-15711') | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), 'https://www.real.discount/offer/web-applications-step-by-step-guide-part-3-27522') | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), 'https://www.real.discount/offer/web-application-step-by-step-guide-27524') | Out-Null
$ws.Hyperlinks.Add($ws.Range("C5"), 'https://www.real.discount/offer/web-applications-step-by-step-guide-part-2-27523') | Out-Null
$ws.Hyperlinks.Add($ws.Range("C6"), 'https://www.real.discount/offer/web-applications-step-by-step-guide-part-4-27521') | Out-Null
$ws.Hyperlinks.Add($ws.Range("C7"), 'https://www.real.discount/offer/18-crucial-cyber-security-tips-29894') | Out-Null
$ws.Hyperlinks.Add($ws.Range("C8"), 'https://www.real.discount/offer/linux-mastery-cli-kali-commands-practice-tests-2024-pro-38400') | Out-Null
$ws.Hyperlinks.Add($ws.Range("C9"), 'https://www.real.discount/offer/exam-ms-900-microsoft-365-fundamentals-mock-exams-35283') | Out-Null
$ws.Hyperlinks.Add($ws.Range("C10"), 'https://www.real.discount/offer/learn-embarcadero-borland-c-builder-in-1-hour-2') | Out-Null
$ws.Hyperlinks.Add($ws.Range("C11"), 'https://www.real.discount/offer/comptia-security-sy0-701-practice-tests-36461') | Out-Null
$ws.Hyperlinks.Add($ws.Range("C12"), 'https://www.real.discount/offer/isc-2-certified-in-cybersecurity-cc-practice-exams-36420') | Out-Null
$ws.Hyperlinks.Add($ws.Range("C13"), 'https://www.real.discount/offer/isc-2-certified-in-cybersecurity-cc-practice-exams-set-2-36444') | Out-Null
$ws.Hyperlinks.Add($ws.Range("C14"), 'https://www.real.discount/offer/pcep-30-02-practice-exams-36925') | Out-Null
$ws.Hyperlinks.Add($ws.Range("C15"), 'https://www.real.discount/offer/css-and-javascript-crash-course-13022') | Out-Null
$ws.Hyperlinks.Add($ws.Range("C16"), 'https://www.real.discount/offer/eccouncil-certified-cybersecurity-technician-39209') | Out-Null
$ws.Hyperlinks.Add($ws.Range("C17"), 'https://www.real.discount/offer/web-automation-and-scraping-using-python-36771') | Out-Null
$ws.Hyperlinks.Add($ws.Range("C18"), 'https://www.real.discount/offer/the-best-chatgpt-ai-course-make-money-with-ai-35563') | Out-Null
$ws.Hyperlinks.Add($ws.Range("C19"), 'https://www.real.discount/offer/learn-azure-bicep-20512') | Out-Null
$ws.Hyperlinks.Add($ws.Range("C20"), 'https://www.real.discount/offer/midjourney-for-beginners-embark-on-your-artistic-journey-35457') | Out-Null
$ws.Hyperlinks.Add($ws.Range("C21"), 'https://www.real.discount/offer/google-forms-o-formularios-de-cero-a-avanzado-32061') | Out-Null

# Restore the Hyperlink cell style that Hyperlinks.Add overrides
$ws.Range("C2:C21").Style = "Hyperlink"
